$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "axios 注销请求"
$ws.Range("C11").Value = "Setting"
$ws.Range("D11").Value = "未做"

$ws.Range("D15").Select()
